# Update TPM-derived NATMI ligand-receptor metrics (Sema3b-Nrp1) for rows 2-17
# (columns G,H,I,J = ligand expression/specificity; M,N,O,P = receptor
# expression/specificity; Q,R,S,T = edge weight/specificity), reflecting
# the refreshed TPM values used to recompute the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 3.883322333333334
$ws.Range("H2").Value = 11.649967
$ws.Range("I2").Value = 0.2039370428985499
$ws.Range("J2").Value = 0.2039370428985498
$ws.Range("M2").Value = 127.3992563333333
$ws.Range("N2").Value = 382.197769
$ws.Range("O2").Value = 0.4838549810199306
$ws.Range("P2").Value = 0.4838549810199307
$ws.Range("Q2").Value = 494.7323773692914
$ws.Range("R2").Value = 4452.591396323623
$ws.Range("S2").Value = 0.09867595402093862
$ws.Range("T2").Value = 0.09867595402093864
$ws.Range("G3").Value = 3.883322333333334
$ws.Range("H3").Value = 11.649967
$ws.Range("I3").Value = 0.2039370428985499
$ws.Range("J3").Value = 0.2039370428985498
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("N3").Value = 178.097596
$ws.Range("O3").Value = 0.2254681108101269
$ws.Range("P3").Value = 0.2254681108101269
$ws.Range("Q3").Value = 230.5367906865924
$ws.Range("R3").Value = 2074.831116179332
$ws.Range("S3").Value = 0.04598129978653984
$ws.Range("T3").Value = 0.04598129978653984
$ws.Range("G4").Value = 3.883322333333334
$ws.Range("H4").Value = 11.649967
$ws.Range("I4").Value = 0.2039370428985499
$ws.Range("J4").Value = 0.2039370428985498
$ws.Range("M4").Value = 16.63275166666667
$ws.Range("N4").Value = 49.898255
$ws.Range("O4").Value = 0.06317022542837675
$ws.Range("P4").Value = 0.06317022542837675
$ws.Range("Q4").Value = 64.5903360119539
$ws.Range("R4").Value = 581.313024107585
$ws.Range("S4").Value = 0.01288274897309793
$ws.Range("T4").Value = 0.01288274897309793
$ws.Range("G5").Value = 3.883322333333334
$ws.Range("H5").Value = 11.649967
$ws.Range("I5").Value = 0.2039370428985499
$ws.Range("J5").Value = 0.2039370428985498
$ws.Range("M5").Value = 59.90262233333334
$ws.Range("N5").Value = 179.707867
$ws.Range("O5").Value = 0.2275066827415657
$ws.Range("P5").Value = 0.2275066827415658
$ws.Range("Q5").Value = 232.6211911322655
$ws.Range("R5").Value = 2093.590720190389
$ws.Range("S5").Value = 0.04639704011797346
$ws.Range("T5").Value = 0.04639704011797347
$ws.Range("G6").Value = 9.654910333333332
$ws.Range("H6").Value = 28.964731
$ws.Range("I6").Value = 0.5070384824688307
$ws.Range("J6").Value = 0.5070384824688307
$ws.Range("M6").Value = 127.3992563333333
$ws.Range("N6").Value = 382.197769
$ws.Range("O6").Value = 0.4838549810199306
$ws.Range("P6").Value = 0.4838549810199307
$ws.Range("Q6").Value = 1230.028396431682
$ws.Range("R6").Value = 11070.25556788514
$ws.Range("S6").Value = 0.2453330953113305
$ws.Range("T6").Value = 0.2453330953113305
$ws.Range("G7").Value = 9.654910333333332
$ws.Range("H7").Value = 28.964731
$ws.Range("I7").Value = 0.5070384824688307
$ws.Range("J7").Value = 0.5070384824688307
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("N7").Value = 178.097596
$ws.Range("O7").Value = 0.2254681108101269
$ws.Range("P7").Value = 0.2254681108101269
$ws.Range("Q7").Value = 573.1721066540749
$ws.Range("R7").Value = 5158.548959886675
$ws.Range("S7").Value = 0.1143210087502809
$ws.Range("T7").Value = 0.1143210087502809
$ws.Range("G8").Value = 9.654910333333332
$ws.Range("H8").Value = 28.964731
$ws.Range("I8").Value = 0.5070384824688307
$ws.Range("J8").Value = 0.5070384824688307
$ws.Range("M8").Value = 16.63275166666667
$ws.Range("N8").Value = 49.898255
$ws.Range("O8").Value = 0.06317022542837675
$ws.Range("P8").Value = 0.06317022542837675
$ws.Range("Q8").Value = 160.5877259382672
$ws.Range("R8").Value = 1445.289533444405
$ws.Range("S8").Value = 0.03202973523841808
$ws.Range("T8").Value = 0.03202973523841808
$ws.Range("G9").Value = 9.654910333333332
$ws.Range("H9").Value = 28.964731
$ws.Range("I9").Value = 0.5070384824688307
$ws.Range("J9").Value = 0.5070384824688307
$ws.Range("M9").Value = 59.90262233333334
$ws.Range("N9").Value = 179.707867
$ws.Range("O9").Value = 0.2275066827415657
$ws.Range("P9").Value = 0.2275066827415658
$ws.Range("Q9").Value = 578.3544473598641
$ws.Range("R9").Value = 5205.190026238777
$ws.Range("S9").Value = 0.1153546431688012
$ws.Range("T9").Value = 0.1153546431688012
$ws.Range("G10").Value = 4.652793333333332
$ws.Range("H10").Value = 13.95838
$ws.Range("I10").Value = 0.2443466784802274
$ws.Range("J10").Value = 0.2443466784802274
$ws.Range("M10").Value = 127.3992563333333
$ws.Range("N10").Value = 382.197769
$ws.Range("O10").Value = 0.4838549810199306
$ws.Range("P10").Value = 0.4838549810199307
$ws.Range("Q10").Value = 592.7624105393577
$ws.Range("R10").Value = 5334.861694854219
$ws.Range("S10").Value = 0.1182283574783335
$ws.Range("T10").Value = 0.1182283574783336
$ws.Range("G11").Value = 4.652793333333332
$ws.Range("H11").Value = 13.95838
$ws.Range("I11").Value = 0.2443466784802274
$ws.Range("J11").Value = 0.2443466784802274
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("N11").Value = 178.097596
$ws.Range("O11").Value = 0.2254681108101269
$ws.Range("P11").Value = 0.2254681108101269
$ws.Range("Q11").Value = 276.2171024504977
$ws.Range("R11").Value = 2485.953922054479
$ws.Range("S11").Value = 0.05509238397966636
$ws.Range("T11").Value = 0.05509238397966637
$ws.Range("G12").Value = 4.652793333333332
$ws.Range("H12").Value = 13.95838
$ws.Range("I12").Value = 0.2443466784802274
$ws.Range("J12").Value = 0.2443466784802274
$ws.Range("M12").Value = 16.63275166666667
$ws.Range("N12").Value = 49.898255
$ws.Range("O12").Value = 0.06317022542837675
$ws.Range("P12").Value = 0.06317022542837675
$ws.Range("Q12").Value = 77.38875606965554
$ws.Range("R12").Value = 696.4988046268999
$ws.Range("S12").Value = 0.01543543476227106
$ws.Range("T12").Value = 0.01543543476227106
$ws.Range("G13").Value = 4.652793333333332
$ws.Range("H13").Value = 13.95838
$ws.Range("I13").Value = 0.2443466784802274
$ws.Range("J13").Value = 0.2443466784802274
$ws.Range("M13").Value = 59.90262233333334
$ws.Range("N13").Value = 179.707867
$ws.Range("O13").Value = 0.2275066827415657
$ws.Range("P13").Value = 0.2275066827415658
$ws.Range("Q13").Value = 278.7145218417178
$ws.Range("R13").Value = 2508.43069657546
$ws.Range("S13").Value = 0.05559050225995647
$ws.Range("T13").Value = 0.05559050225995648
$ws.Range("G14").Value = 0.8507443333333334
$ws.Range("H14").Value = 2.552233
$ws.Range("I14").Value = 0.04467779615239207
$ws.Range("J14").Value = 0.04467779615239207
$ws.Range("M14").Value = 127.3992563333333
$ws.Range("N14").Value = 382.197769
$ws.Range("O14").Value = 0.4838549810199306
$ws.Range("P14").Value = 0.4838549810199307
$ws.Range("Q14").Value = 108.3841953964641
$ws.Range("R14").Value = 975.457758568177
$ws.Range("S14").Value = 0.021617574209328
$ws.Range("T14").Value = 0.021617574209328
$ws.Range("G15").Value = 0.8507443333333334
$ws.Range("H15").Value = 2.552233
$ws.Range("I15").Value = 0.04467779615239207
$ws.Range("J15").Value = 0.04467779615239207
$ws.Range("M15").Value = 59.36586533333332
$ws.Range("N15").Value = 178.097596
$ws.Range("O15").Value = 0.2254681108101269
$ws.Range("P15").Value = 0.2254681108101269
$ws.Range("Q15").Value = 50.50517352576311
$ws.Range("R15").Value = 454.546561731868
$ws.Range("S15").Value = 0.01007341829363979
$ws.Range("T15").Value = 0.0100734182936398
$ws.Range("G16").Value = 0.8507443333333334
$ws.Range("H16").Value = 2.552233
$ws.Range("I16").Value = 0.04467779615239207
$ws.Range("J16").Value = 0.04467779615239207
$ws.Range("M16").Value = 16.63275166666667
$ws.Range("N16").Value = 49.898255
$ws.Range("O16").Value = 0.06317022542837675
$ws.Range("P16").Value = 0.06317022542837675
$ws.Range("Q16").Value = 14.15021922815722
$ws.Range("R16").Value = 127.351973053415
$ws.Range("S16").Value = 0.00282230645458967
$ws.Range("T16").Value = 0.00282230645458967
$ws.Range("G17").Value = 0.8507443333333334
$ws.Range("H17").Value = 2.552233
$ws.Range("I17").Value = 0.04467779615239207
$ws.Range("J17").Value = 0.04467779615239207
$ws.Range("M17").Value = 59.90262233333334
$ws.Range("N17").Value = 179.707867
$ws.Range("O17").Value = 0.2275066827415657
$ws.Range("P17").Value = 0.2275066827415658
$ws.Range("Q17").Value = 50.96181650189013
$ws.Range("R17").Value = 458.6563485170111
$ws.Range("S17").Value = 0.01016449719483461
$ws.Range("T17").Value = 0.01016449719483461
